$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 170911
$ws.Range("C4").Value = 161725
$ws.Range("C5").Value = 9186
$ws.Range("C8").Value = 65.81999999999999
